# Add a new "RiparianBufferArea" column (K) to the PlanningUnits sheet,
# populate it with a per-row formula against the existing ncols column (J),
# and mark the Identifier column (A) as Text-formatted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PlanningUnits")

# Header for the new column.
$ws.Range("K1").Value = "RiparianBufferArea"

# Per-row formula: half of column J's value, row by row.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 11).Formula = '=$J' + $r + '/2'
}

# Give the new column a sensible width (mirrors the other bestFit columns).
$ws.Columns.Item(11).ColumnWidth = 17.5

# Identifier column becomes Text-formatted.
$ws.Range("A1:A6").NumberFormat = "@"

# Leave the authored selection on B1.
$ws.Range("B1").Select() | Out-Null

# Sheet is printed in portrait orientation.
$ws.PageSetup.Orientation = 1
